$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (SamsungElec / 005930.KS) - only the final score (N) changes
$ws.Range("N2").Value = 85.8724807945396

# Row 3 now represents 240810.KS (ticker string moved here)
$ws.Range("B3").Value = "240810.KS,0P00017YB3,330568"
$ws.Range("C3").Value = "240810.KS"
$ws.Range("D3").Value = 61300
$ws.Range("F3").Value = 8.109999999999999
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 63
$ws.Range("I3").Value = 70
$ws.Range("J3").Value = 66
$ws.Range("K3").Value = 59.8
$ws.Range("N3").Value = 85.8724807945396

# Row 4 now represents 058470.KS (ticker string moved here)
$ws.Range("B4").Value = "058470.KS,0P0000ASU1,98886"
$ws.Range("C4").Value = "058470.KS"
$ws.Range("D4").Value = 68300
$ws.Range("F4").Value = 25.55
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 36
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 63
$ws.Range("K4").Value = 57.8
$ws.Range("N4").Value = 85.8724807945396

# Row 5 now represents DB HiTek / 000990.KS
$ws.Range("B5").Value = "DB HiTek"
$ws.Range("C5").Value = "000990.KS"
$ws.Range("D5").Value = 63600
$ws.Range("F5").Value = 2.75
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 36
$ws.Range("I5").Value = 56
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 57.2
$ws.Range("N5").Value = 85.8724807945396

# Row 6 now represents SK hynix / 000660.KS
$ws.Range("B6").Value = "SK hynix"
$ws.Range("C6").Value = "000660.KS"
$ws.Range("D6").Value = 530000
$ws.Range("F6").Value = 1.8
$ws.Range("G6").Value = 20
$ws.Range("H6").Value = 60
$ws.Range("I6").Value = 60
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 55.8
$ws.Range("N6").Value = 85.8724807945396

# Row 7 (403870.KS) keeps its identity, but data values refreshed
$ws.Range("D7").Value = 30250
$ws.Range("F7").Value = 6.7
$ws.Range("G7").Value = 20
$ws.Range("H7").Value = 53
$ws.Range("I7").Value = 46
$ws.Range("J7").Value = 56
$ws.Range("K7").Value = 50.2
$ws.Range("N7").Value = 85.8724807945396
